$wb = $excel.ActiveWorkbook

# --- Sheet1: CalculationAmountAndPrice ---
$ws1 = $wb.Worksheets.Item("CalculationAmountAndPrice")

# Row 2 switches to the 18-decimals SSRT token, so the traded amount (D2)
# is expressed with fewer "whole" units than before.
$ws1.Range("C2").Value = "Snip 20 - SSRT ( 18 decimals)"
$ws1.Range("D2").Value = 100

# Fix the sell-side price formulas on rows 4 and 5 (they were computing the
# inverse of the correct price).
$ws1.Range("J4").Formula = "=1/ (F4/D4)"
$ws1.Range("J5").Formula = "=1 / (F5/D5)"

# Column C needs to be a bit wider to fit the new label.
$ws1.Columns.Item(3).ColumnWidth = 26.1666666666667

# --- Sheet2: CalculationLiquidity ---
$ws2 = $wb.Worksheets.Item("CalculationLiquidity")
$ws2.Columns.Item(13).ColumnWidth = 14.665

# Update the view selections (sheet2 first, then sheet1 last so sheet1
# remains the active/tab-selected sheet).
$ws2.Activate()
$ws2.Range("K5").Select()

$ws1.Activate()
$ws1.Application.ActiveWindow.Zoom = 100
$ws1.Range("J6").Select()

$wb.Save()
